# Auto-generated Excel COM-interop script
# Updates symbol list (coin prices/volumes, and row reshuffle for B6:E18)
# per commit: "Updated symbol list on Wed Feb  1 08:56:30 UTC 2023 with GitHub Actions"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume(1h)) hold text that looks numeric/percent-like
# (e.g. "307.69", "-0.91%"). Force Text number format before assignment so Excel
# keeps them as literal strings instead of coercing to number/percentage types,
# matching the original inlineStr/text storage.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "307.69"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-0.91%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "37.40"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-0.40%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.102"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "0.35%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07790"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "0.46%"
$ws.Range("B6").Value = "KuCoinToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "8.248"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "0.35%"
$ws.Range("B7").Value = "FTXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.880"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-0.57%"
$ws.Range("B8").Value = "BTSEToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "2.988"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "3.51%"
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9289"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "0.89%"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1089"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-10.13%"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1916"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-0.21%"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08906"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-4.28%"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03333"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-2.89%"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09592"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.88%"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001387"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.57%"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005735"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-2.11%"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.543"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.23%"
$ws.Range("B18").Value = "GateToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.436"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "1.76%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3475"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.285"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "19.21%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1274"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-0.66%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2507"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-3.40%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04385"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "0.58%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001202"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-1.05%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004250"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-0.03%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001303"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "0.21%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02164"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "3.66%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05023"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "0.96%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007471"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-2.50%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1347"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "0.18%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.008661"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-12.01%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002115"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "6.20%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.008007"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-9.92%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006563"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-1.67%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000752"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.08%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.002877"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-1.95%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.001003"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-16.60%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002106"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.08%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002006"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.08%"
